# "New air dehumidification set" -- adds a Multifamily Housing row replacement
# for the former "Apartment Building" entry, fills in the newly-introduced
# "Unoccupied Turndown" (H) / "Moisture control" (F) style placeholder cells
# with an em-dash ("—") where data is not available/applicable, and moves the
# active selection to the new F12:F15 block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dash = [char]0x2014

# Rename "Apartment Building" -> "Multifamily Housing"
$ws.Range("A8").Value = "Multifamily Housing"

# Column H ("Unoccupied Turndown") now carries a placeholder "-" for every
# row that previously had no entry for it.
foreach ($r in 2..10) {
    $ws.Range("H$r").Value = $dash
}

# Rows 4-9 pick up placeholder dashes across the previously-empty columns.
$ws.Range("D4").Value = $dash
$ws.Range("E4").Value = $dash
$ws.Range("F4").Value = $dash
$ws.Range("G4").Value = $dash
$ws.Range("G4").Style = "Normal"

foreach ($r in 5..8) {
    $ws.Range("B$r").Value = $dash
    $ws.Range("D$r").Value = $dash
    $ws.Range("E$r").Value = $dash
    $ws.Range("F$r").Value = $dash
    $ws.Range("G$r").Value = $dash
    $ws.Range("G$r").Style = "Normal"
}

$ws.Range("D9").Value = $dash
$ws.Range("E9").Value = $dash
$ws.Range("F9").Value = $dash
$ws.Range("G9").Value = $dash
$ws.Range("G9").Style = "Normal"

# Newly introduced "Moisture control" (F) placeholder dashes for the
# hospital-room sub-rows.
foreach ($r in 12..15) {
    $ws.Range("F$r").Value = $dash
}

# Move the active selection to the new F12:F15 block, matching the
# post-edit sheet view.
$ws.Range("F12:F15").Select()
